$d = $word.ActiveDocument

# 1) Rename the first heading from "Objet positionUpdate" to "Objet geoPos"
$d.Content.Find.Execute("Objet positionUpdate", $true, $false, $false, $false, $false, $true, 1, $false, "Objet geoPos", 2)

# 2) Make room right after the renamed heading paragraph by inserting a
#    fresh empty paragraph, then replace that empty paragraph with the
#    new "positionUpdate" summary table plus the "Type position" heading
#    (InsertXML on a collapsed range absorbs the paragraph it sits in).
$p1 = $d.Paragraphs(1)
$rng = $p1.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$rng2 = $p2.Range
$rng2.Collapse(0)
$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblStyle w:val="MediumShading1-Accent1"/><w:tblW w:type="auto" w:w="0"/><w:tblLayout w:type="fixed"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="2040"/><w:gridCol w:w="2040"/><w:gridCol w:w="2040"/><w:gridCol w:w="2040"/><w:gridCol w:w="2040"/><w:gridCol w:w="2040"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1701"/></w:tcPr><w:p><w:r><w:t>Nom de balise</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1984"/></w:tcPr><w:p><w:r><w:t>Champ correspondant</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1134"/></w:tcPr><w:p><w:r><w:t>Format</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1417"/></w:tcPr><w:p><w:r><w:t>Cardinalité</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4535"/></w:tcPr><w:p><w:r><w:t>Description</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1701"/></w:tcPr><w:p><w:r><w:t>Exemple</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1701"/></w:tcPr><w:p><w:r><w:t>positionUpdate</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1984"/></w:tcPr><w:p><w:r><w:t>Position</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1134"/></w:tcPr><w:p><w:r><w:t>cf. type position</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1417"/></w:tcPr><w:p><w:r><w:t>0..n</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4535"/></w:tcPr><w:p><w:r><w:t>Objet de détail de la position de chaque ressource</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1701"/></w:tcPr><w:p><w:r/></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Type position</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@
$rng2.InsertXML($xmlFrag)
